$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fix: normalize phone-number formatting to "(xxx) xxx-xxxx" style ---
$ws.Range("D8").Value  = "(416) 508-1780"
$ws.Range("D13").Value = "(778) 709-1769"
$ws.Range("D15").Value = "(778) 846-9554"
$ws.Range("D16").Value = "(250) 386-5311-2030"
$ws.Range("D17").Value = "(250) 386-5311-2030"

# --- View state: zoom in and move the active selection down past the table ---
$null = $ws.Activate()
$excel.ActiveWindow.Zoom = 150
$null = $ws.Range("D18").Select()

# --- Drop the leftover explicit row heights from the recovered file so rows
#     fall back to the sheet's default height again ---
$null = $ws.Range("A1:D17").EntireRow.AutoFit()
